$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAM Variable Changes")

# Data for new "Deleted variable" rows (18 new rows: 69-86)
# Columns: A=Type, B=Variable Type, C=Old Name, E=Input Page, F=Default/Reason, G=Handled in Version Upgrader, H=Ty
$rows = @(
    @{ C = "m_dot_htf_ref"; E = "Molten Salt Tower Power Block" },
    @{ C = "T_pb_out"; E = "Molten Salt Tower Power Block" },
    @{ C = "mode"; E = "Molten Salt Tower Power Block" },
    @{ C = "fthr_ok"; E = "Molten Salt Tower Power Block" },
    @{ C = "pb_fixed_par_cntl"; E = "Molten Salt Tower Power Block" },
    @{ C = "dt_cold"; E = "Molten Salt Tower Power Block" },
    @{ C = "dt_hot"; E = "Molten Salt Tower Power Block" },
    @{ C = "hx_config"; E = "Molten Salt Tower Power Block" },
    @{ C = "is_hx"; E = "Molten Salt Tower Power Block" },
    @{ C = "tech_type"; E = "Molten Salt Tower Power Block" },
    @{ C = "deg_wind"; E = "Molten Salt Tower Receiver" },
    @{ C = "P_htf"; E = "Molten Salt Tower Receiver" },
    @{ C = "T_salt_cold"; E = "Molten Salt Power Block" },
    @{ C = "HTF"; E = "Molten Salt Power Block" },
    @{ C = "Design_power"; E = "Molten Salt Power Block" },
    @{ C = "csp.pt.pwrb.min_restart_time"; E = "Molten Salt Power Block" },
    @{ C = "csp.pt.rec.max_rec_flux"; E = "Molten Salt Power Block" },
    @{ C = "store_fluid"; E = "Molten Salt Tower Storage" }
)

$startRow = 69
$templateRow = 68
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Copy formatting (A:H only) from an existing data row down to the new row
    $ws.Range("A$($templateRow):H$($templateRow)").Copy()
    $ws.Range("A$($r):H$($r)").PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = "Deleted variable"
    $ws.Cells.Item($r, 2).Value = "number"
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = "not used"
    $ws.Cells.Item($r, 7).Value = "N/A"
    $ws.Cells.Item($r, 8).Value = "Ty"
}
$excel.CutCopyMode = 0

# Widen column C slightly, as seen in the diff
$ws.Columns.Item(3).ColumnWidth = 28

# Update the data validation range for column A to cover the new rows
$ws.Range("A2:A86").Validation.Delete()
$ws.Range("A2:A86").Validation.Add(3, 1, 1, "Types")
$ws.Range("A2:A86").Validation.IgnoreBlank = $true
$ws.Range("A2:A86").Validation.InCellDropdown = $true
$ws.Range("A2:A86").Validation.ShowInput = $true
$ws.Range("A2:A86").Validation.ShowError = $true

# Update the active view to match the diff
$ws.Application.ActiveWindow.ScrollRow = 52
$ws.Range("C77").Select()
